# Add newly-collected studies (rows 96-103) to the literature table, plus
# two trailing blank (but formatted) rows, mirroring the author's "Add files
# via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Em dash used in the "Big data analytics architecture design—an application..."
# title; built from its code point so the .ps1 source stays plain ASCII.
$emDash = [char]0x2014

# New rows of source data: Title / Retrieval source / Type / Year / Search method
$newRows = @(
    @("Euromicro Conference on Software Engineering and Advanced Applications", "Flexible System-Level Monitoring of Heterogeneous Big Data Streaming Systems", "Conference", 2018, "Automatic search"),
    @("IEEE Conference on Software Quality, Reliability and Security Companion", "Defect Prediction Based on the Characteristics of Multilayer Structure of Software Network", "Conference", 2018, " Manual search"),
    @("Computers and Industrial Engineering", ("Big data analytics architecture design" + $emDash + "an application in manufacturing systems"), "Journal", 2019, "Automatic search"),
    @("IEEE Access", "Big Data Opportunities: System Health Monitoring and Management", "Journal", 2019, "Automatic search"),
    @("Future Generation Computer Systems", "BDWatchdog: Real-time monitoring and profiling of Big Data applications and frameworks", "Journal", 2018, " Manual search"),
    @("ACM/SIGAPP Symposium on Applied Computing", "Evaluation of ACE properties of traditional SQL and NoSQL big data systems", "Conference", 2019, " Manual search"),
    @("IEEE International Conference on Control Systems and Computer Science", "Systems Monitoring and Big Data Analysis Using the Elasticsearch System", "Conference", 2019, "Automatic search"),
    @("IEEE Transactions on Information Forensics and Security", "Collaboration- and Fairness-Aware Big Data Management in Distributed Clouds", "Journal", 2019, " Manual search")
)

$startRow = 96
$lastDataRow = $startRow + $newRows.Count - 1   # 103
$blankRow1 = $lastDataRow + 1                   # 104
$blankRow2 = $lastDataRow + 2                   # 105

# Carry the last existing data row's formatting down across the new data
# rows (A95:E95 -> A96:E103) before writing values, so the appended rows
# look like the rest of the table.
$ws.Range("A95:E95").Copy()
$ws.Range("A96:E" + $lastDataRow).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# Rows 104 and 105 stay blank, with only columns A/B formatted (matching the
# trailing empty rows left under the table).
$ws.Range("A95:B95").Copy()
$ws.Range("A" + $blankRow1 + ":B" + $blankRow1).PasteSpecial(-4122)
$ws.Range("A" + $blankRow2 + ":B" + $blankRow2).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# A handful of existing rows (26, 58, 74, 84) carried a stray duplicate
# "Times New Roman 8, explicit black" font instead of the normal body style;
# normalize them to match the rest of column A/B.
$normalizeRows = @(26, 58, 74, 84)
foreach ($rowNum in $normalizeRows) {
    $ws.Range("A2:B2").Copy()
    $ws.Range("A" + $rowNum + ":B" + $rowNum).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

$ws.Range("A103").Select()
